$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 updates
$ws.Range("B12").Value = 85
$ws.Range("C12").Value = -3.6
$ws.Range("E12").Value = "81.4/140"
